# Prox Ops Definition Illustration PP.pptx
#
# 1) Update the cached "datetimeFigureOut" field text from 3/10/2014 to
#    3/21/2014 everywhere it appears (the slide master and all 11 slide
#    layouts each carry one such field in their Date placeholder).
# 2) On Slide 3, re-label the three coordinate-axis call-outs, changing
#    the wording inside the smart-quoted text boxes:
#       "In Track"    -> "In-Track"
#       "Out of Plane"-> "Out-of-Plane"
#       "Out of Track"-> "Cross-Track"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder refresh (Slide Master + every Custom Layout)
# ---------------------------------------------------------------------
$oldDate = "3/10/2014"
$newDate = "3/21/2014"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 3 axis-label call-outs
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# "In Track" -> "In-Track"  (quotes untouched, word inside re-typed)
$shInTrack = $s3.Shapes.Item("TextBox 22")
$paraInTrack = $shInTrack.TextFrame.TextRange.Paragraphs(2)
$paraInTrack.Characters(2, 8).Text = "In-Track"

# "Out of Plane" -> "Out-of-Plane"
$shOutOfPlane = $s3.Shapes.Item("TextBox 23")
$paraOutOfPlane = $shOutOfPlane.TextFrame.TextRange.Paragraphs(2)
$paraOutOfPlane.Characters(2, 12).Text = "Out-of-Plane"

# "Out of Track" -> "Cross-Track" (leading curly quote retyped together
# with the new word, trailing curly quote left as its own run)
$shOutOfTrack = $s3.Shapes.Item("TextBox 24")
$paraOutOfTrack = $shOutOfTrack.TextFrame.TextRange.Paragraphs(2)
$paraOutOfTrack.Characters(1, 13).Text = [char]8220 + "Cross-Track"
